$d = $word.ActiveDocument

# Explicitly set PageBreakBefore = False on every paragraph in the body.
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = 0
}

# Also stamp the same explicit (default) value onto the paragraph styles
# (Heading 1-6, Title, Subtitle) used by/available to this document.
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", `
                "Heading 5", "Heading 6", "Title", "Subtitle")
foreach ($name in $styleNames) {
    $s = $d.Styles($name)
    $s.ParagraphFormat.PageBreakBefore = 0
}
